$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# Copy the date cell from the last existing ballot row so the new row's
# AM cell picks up the same date number format (style) without creating a
# brand new custom number format in styles.xml.
$ws.Range("AM23").Copy($ws.Range("AM24"))

# New ballot row for René Cárdenas
$ws.Range("A24").Value = "René Cárdenas"
$ws.Range("B24").Value = "x"
$ws.Range("D24").Value = "x"
$ws.Range("H24").Value = "x"
$ws.Range("I24").Value = "x"
$ws.Range("L24").Value = "x"
$ws.Range("O24").Value = "x"
$ws.Range("Q24").Value = "x"
$ws.Range("T24").Value = "x"
$ws.Range("U24").Value = "x"
$ws.Range("V24").Value = "x"
$ws.Range("AK24").Value = 10
$ws.Range("AL24").Value = "Twitter"
$ws.Range("AM24").Value = 43441

# Update the saved selection to match the new active cell
$ws.Activate() | Out-Null
$ws.Range("A24").Select() | Out-Null
